$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 505
$ws.Range("J38").Value = 2222
$ws.Range("L38").Value = 6666
$ws.Range("N38").Value = -7410

$ws.Range("H41").Value = 229.2
$ws.Range("I41").Value = 211.75
$ws.Range("J41").Value = 299
$ws.Range("K41").Value = 211.75
$ws.Range("L41").Value = 299
$ws.Range("M41").Value = 228.25
$ws.Range("N41").Value = -1179

$ws.Range("H70").Value = 1916.4445
$ws.Range("I70").Value = 1639.25
$ws.Range("K70").Value = 4917.75
$ws.Range("M70").Value = -4647.75

$ws.Range("H73").Value = 1916.4445
$ws.Range("I73").Value = 1639.25
$ws.Range("K73").Value = 4917.75
$ws.Range("M73").Value = -3981.75

$ws.Range("H76").Value = 90917640
$ws.Range("J76").Value = 9894.5
$ws.Range("L76").Value = 9894.5
$ws.Range("N76").Value = -10524.5

$ws.Range("H79").Value = 90917640
$ws.Range("J79").Value = 9894.5
$ws.Range("L79").Value = 9894.5
$ws.Range("N79").Value = -12078.5

$ws.Range("H82").Value = 348
$ws.Range("I82").Value = 348
$ws.Range("K82").Value = 1044
$ws.Range("M82").Value = -638

$ws.Range("H85").Value = 348
$ws.Range("I85").Value = 348
$ws.Range("K85").Value = 1044
$ws.Range("M85").Value = 360

$ws.Range("H92").Value = 734.8333
$ws.Range("I92").Value = 682.1
$ws.Range("J92").Value = 998.5
$ws.Range("K92").Value = 682.1
$ws.Range("L92").Value = 998.5
$ws.Range("M92").Value = 565.9
$ws.Range("N92").Value = -3494.5

$ws.Range("H96").Value = 1765.75
$ws.Range("I96").Value = 1332.8334
$ws.Range("J96").Value = 3064.5
$ws.Range("K96").Value = 3998.5002
$ws.Range("L96").Value = 9193.5
$ws.Range("M96").Value = -2625.5002
$ws.Range("N96").Value = -11939.5

$ws.Range("H99").Value = 9831.333000000001
$ws.Range("I99").Value = 997
$ws.Range("J99").Value = 27500
$ws.Range("K99").Value = 2991
$ws.Range("L99").Value = 82500
$ws.Range("M99").Value = -1493
$ws.Range("N99").Value = -85496

$ws.Range("H141").Value = 3616.476
$ws.Range("I141").Value = 3135.9443
$ws.Range("K141").Value = 9407.832900000001
$ws.Range("M141").Value = -4227.832900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4231.4443
$ws.Range("J63").Value = 5706.273
$ws.Range("L63").Value = 5706.273
$ws.Range("N63").Value = -7078.273

$ws.Range("H66").Value = 4231.4443
$ws.Range("J66").Value = 5706.273
$ws.Range("L66").Value = 28531.365
$ws.Range("N66").Value = -35395.36500000001

$ws.Range("H102").Value = 3448.4783
$ws.Range("I102").Value = 3464.318
$ws.Range("J102").Value = 3100
$ws.Range("K102").Value = 3464.318
$ws.Range("L102").Value = 3100
$ws.Range("M102").Value = -1842.318
$ws.Range("N102").Value = -6344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 66682932
$ws.Range("I58").Value = 90921500
$ws.Range("J58").Value = 26871.25
$ws.Range("K58").Value = 90921500
$ws.Range("L58").Value = 26871.25
$ws.Range("M58").Value = -90921297
$ws.Range("N58").Value = -27277.25

$ws.Range("H122").Value = 3481.973
$ws.Range("I122").Value = 3364.2
$ws.Range("J122").Value = 3727.3333
$ws.Range("K122").Value = 10092.6
$ws.Range("L122").Value = 11181.9999
$ws.Range("M122").Value = -7642.599999999999
$ws.Range("N122").Value = -16081.9999

$ws.Range("H132").Value = 5328.697
$ws.Range("I132").Value = 4583.7744
$ws.Range("K132").Value = 13751.3232
$ws.Range("M132").Value = -11221.3232

$ws.Range("H136").Value = 66682932
$ws.Range("I136").Value = 90921500
$ws.Range("J136").Value = 26871.25
$ws.Range("K136").Value = 272764500
$ws.Range("L136").Value = 80613.75
$ws.Range("M136").Value = -272761950
$ws.Range("N136").Value = -85713.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3508.3333
$ws.Range("I34").Value = 50
$ws.Range("K34").Value = 150
$ws.Range("M34").Value = -66

$ws.Range("H50").Value = 4040.8
$ws.Range("I50").Value = 3801
$ws.Range("K50").Value = 11403
$ws.Range("M50").Value = -10922

$ws.Range("H52").Value = 3474501
$ws.Range("J52").Value = 3474501
$ws.Range("L52").Value = 10423503
$ws.Range("N52").Value = -10424035

$ws.Range("H53").Value = 4040.8
$ws.Range("I53").Value = 3801
$ws.Range("K53").Value = 11403
$ws.Range("M53").Value = -10922

$ws.Range("H119").Value = 1968.25
$ws.Range("I119").Value = 1535.1428
$ws.Range("J119").Value = 5000
$ws.Range("K119").Value = 4605.428400000001
$ws.Range("L119").Value = 15000
$ws.Range("M119").Value = 232.5715999999993
$ws.Range("N119").Value = -24676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11444286
$ws.Range("I11").Value = 10020000
$ws.Range("J11").Value = 15005000
$ws.Range("K11").Value = 10020000
$ws.Range("L11").Value = 15005000
$ws.Range("M11").Value = -10019861
$ws.Range("N11").Value = -15005278

$ws.Range("H107").Value = 1346.0834
$ws.Range("I107").Value = 753.8570999999999
$ws.Range("J107").Value = 2175.2
$ws.Range("K107").Value = 753.8570999999999
$ws.Range("L107").Value = 2175.2
$ws.Range("M107").Value = 1166.1429
$ws.Range("N107").Value = -6015.2

$ws.Range("H132").Value = 27780622
$ws.Range("I132").Value = 55557736
$ws.Range("K132").Value = 166673208
$ws.Range("M132").Value = -166670678

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 772.25
$ws.Range("I22").Value = 536.8
$ws.Range("K22").Value = 536.8
$ws.Range("M22").Value = -241.8

$ws.Range("H27").Value = 772.25
$ws.Range("I27").Value = 536.8
$ws.Range("K27").Value = 536.8
$ws.Range("M27").Value = -429.8

$ws.Range("H46").Value = 33334338
$ws.Range("I46").Value = 1019.8
$ws.Range("J46").Value = 50000996
$ws.Range("K46").Value = 1019.8
$ws.Range("L46").Value = 50000996
$ws.Range("M46").Value = -831.8
$ws.Range("N46").Value = -50001372

$ws.Range("H132").Value = 3621.45
$ws.Range("J132").Value = 4387.5
$ws.Range("L132").Value = 13162.5
$ws.Range("N132").Value = -18222.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 22099.8
$ws.Range("I62").Value = 18875
$ws.Range("J62").Value = 24249.666
$ws.Range("K62").Value = 18875
$ws.Range("L62").Value = 24249.666
$ws.Range("M62").Value = -18251
$ws.Range("N62").Value = -25497.666

$ws.Range("H65").Value = 22099.8
$ws.Range("I65").Value = 18875
$ws.Range("J65").Value = 24249.666
$ws.Range("K65").Value = 94375
$ws.Range("L65").Value = 121248.33
$ws.Range("M65").Value = -91255
$ws.Range("N65").Value = -127488.33
